$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 560, shifting existing rows 560:589 down to 561:590.
$ws.Rows.Item(560).Insert()

# Populate the newly-inserted row 560 with the new data record.
$ws.Cells.Item(560, 1).Value = 10
$ws.Cells.Item(560, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(560, 3).Value = "La Araucanía"
$ws.Cells.Item(560, 4).Value = 44753
$ws.Cells.Item(560, 5).Value = 9
$ws.Cells.Item(560, 6).Value = 100112045
$ws.Cells.Item(560, 7).Value = "Zapallo"
$ws.Cells.Item(560, 8).Value = "Camote"
$ws.Cells.Item(560, 9).Value = "1a (guarda)"
$ws.Cells.Item(560, 10).Value = 900
$ws.Cells.Item(560, 11).Value = 700
$ws.Cells.Item(560, 12).Value = 700
$ws.Cells.Item(560, 13).Value = 700
$ws.Cells.Item(560, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(560, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(560, 16).Value = 700
$ws.Cells.Item(560, 17).Value = 1
$ws.Cells.Item(560, 18).Value = "Hortaliza"
